# Update cryptocurrency price/volume data per latest GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.063.23'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '3.433.93'
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("D4").Value = "'" + '0.999'
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = "'" + '410.24'
$ws.Range("E5").Value = '  +1.65%  '
$ws.Range("D6").Value = "'" + '129.08'
$ws.Range("E6").Value = '  -2.17%  '
$ws.Range("D7").Value = "'" + '0.630'
$ws.Range("E7").Value = '  +7.73%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = "'" + '0.734'
$ws.Range("E9").Value = '  +10.50%  '
$ws.Range("D10").Value = "'" + '0.147'
$ws.Range("E10").Value = '  +23.86%  '
$ws.Range("D11").Value = "'" + '42.57'
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("D12").Value = "'" + '0.0000223'
$ws.Range("E12").Value = '  +75.73%  '
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").Value = '3.974.71'
$ws.Range("E14").Value = '  +2.00%  '
$ws.Range("E15").Value = '  +7.64%  '
$ws.Range("D16").Value = "'" + '8.93'
$ws.Range("E16").Value = '  +6.79%  '
$ws.Range("D17").Value = '3.478.31'
$ws.Range("E17").Value = '  +3.45%  '
$ws.Range("D18").Value = "'" + '12.38'
$ws.Range("E18").Value = '  +13.89%  '
$ws.Range("E19").Value = '  +6.11%  '
$ws.Range("D20").Value = '62.013.58'
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("D21").Value = "'" + '391.57'
$ws.Range("E21").Value = '  +25.10%  '
$ws.Range("D22").Value = "'" + '90.22'
$ws.Range("E22").Value = '  +7.04%  '
$ws.Range("E23").Value = '  +1.06%  '
$ws.Range("E24").Value = '  +4.51%  '
$ws.Range("D25").Value = "'" + '3.22'
$ws.Range("E25").Value = '  +4.34%  '
$ws.Range("D26").Value = "'" + '32.76'
$ws.Range("E26").Value = '  +11.70%  '
$ws.Range("D27").Value = "'" + '8.77'
$ws.Range("E27").Value = '  +7.23%  '
$ws.Range("D28").Value = "'" + '4.80'
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = "'" + '7.66'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = "'" + '2.77'
$ws.Range("E30").Value = '  +2.35%  '
$ws.Range("D31").Value = "'" + '0.119'
$ws.Range("E31").Value = '  +2.95%  '
$ws.Range("D32").Value = "'" + '11.97'
$ws.Range("E32").Value = '  +6.05%  '
$ws.Range("D33").Value = "'" + '43.99'
$ws.Range("E33").Value = '  +5.74%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").Value = "'" + '0.0502'
$ws.Range("E36").Value = '  +5.40%  '
$ws.Range("E37").Value = '  +3.91%  '
$ws.Range("D38").Value = "'" + '0.997'
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").Value = "'" + '3.40'
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").Value = "'" + '2.92'
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("E41").Value = '  +7.39%  '
$ws.Range("D42").Value = "'" + '0.314'
$ws.Range("E42").Value = '  +8.36%  '
$ws.Range("D43").Value = "'" + '141.52'
$ws.Range("E43").Value = '  +2.14%  '
$ws.Range("D44").Value = "'" + '1.99'
$ws.Range("E44").Value = '  +0.76%  '
$ws.Range("D45").Value = "'" + '4.11'
$ws.Range("E45").Value = '  +3.95%  '
$ws.Range("D46").Value = "'" + '2.44'
$ws.Range("E46").Value = '  +10.17%  '
$ws.Range("D47").Value = "'" + '16.79'
$ws.Range("E47").Value = '  +1.19%  '
$ws.Range("D48").Value = "'" + '21.98'
$ws.Range("E48").Value = '  +3.85%  '
$ws.Range("D49").Value = '2.125.03'
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("B50").Value = 'PEPE'
$ws.Range("C50").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D50").Value = '0.0₃0480'
$ws.Range("E50").Value = '  +65.08%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'" + '0.131'
$ws.Range("E51").Value = '  +16.78%  '
